$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $ref, $val) {
    $rng = $ws.Range($ref)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.ClearFormats()
}

# Row 2
Set-TextValue $ws "D2" '20.539.81'
Set-TextValue $ws "E2" '  -6.94%  '

# Row 3
Set-TextValue $ws "D3" '1.451.51'
Set-TextValue $ws "E3" '  -6.83%  '

# Row 4
Set-TextValue $ws "D4" '1.008'
Set-TextValue $ws "E4" '  +0.84%  '

# Row 5
Set-TextValue $ws "D5" '1.008'
Set-TextValue $ws "E5" '  +0.88%  '

# Row 6
Set-TextValue $ws "D6" '277.28'
Set-TextValue $ws "E6" '  -4.99%  '

# Row 7
Set-TextValue $ws "D7" '0.3718'
Set-TextValue $ws "E7" '  -6.30%  '

# Row 8
Set-TextValue $ws "D8" '0.3105'
Set-TextValue $ws "E8" '  -4.28%  '

# Row 9
Set-TextValue $ws "D9" '41.34'
Set-TextValue $ws "E9" '  -6.75%  '

# Row 10
Set-TextValue $ws "D10" '1.016'
Set-TextValue $ws "E10" '  -6.25%  '

# Row 11
Set-TextValue $ws "D11" '0.06649'
Set-TextValue $ws "E11" '  -8.76%  '

# Row 12
Set-TextValue $ws "D12" '1.007'
Set-TextValue $ws "E12" '  +0.68%  '

# Row 13
Set-TextValue $ws "D13" '5.425'
Set-TextValue $ws "E13" '  -5.21%  '

# Row 14
Set-TextValue $ws "D14" '17.41'
Set-TextValue $ws "E14" '  -7.91%  '

# Row 15
Set-TextValue $ws "B15" 'Chainlink'
Set-TextValue $ws "C15" 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-TextValue $ws "D15" '6.181'
Set-TextValue $ws "E15" '  -7.27%  '

# Row 16
Set-TextValue $ws "B16" 'WrappedEther'
Set-TextValue $ws "C16" 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextValue $ws "D16" '1.456.43'
Set-TextValue $ws "E16" '  -6.49%  '

# Row 17
Set-TextValue $ws "D17" '0.00001027'
Set-TextValue $ws "E17" '  -8.94%  '

# Row 18
Set-TextValue $ws "D18" '0.06077'
Set-TextValue $ws "E18" '  -7.82%  '

# Row 19
Set-TextValue $ws "D19" '77.93'
Set-TextValue $ws "E19" '  -7.22%  '

# Row 20
Set-TextValue $ws "D20" '1.009'
Set-TextValue $ws "E20" '  +1.05%  '

# Row 21
Set-TextValue $ws "D21" '5.752'
Set-TextValue $ws "E21" '  -8.65%  '

# Row 22
Set-TextValue $ws "D22" '14.66'
Set-TextValue $ws "E22" '  -6.27%  '

# Row 23
Set-TextValue $ws "D23" '10.96'
Set-TextValue $ws "E23" '  -3.46%  '

# Row 24
Set-TextValue $ws "D24" '2.315'
Set-TextValue $ws "E24" '  -2.17%  '

# Row 25
Set-TextValue $ws "D25" '20.582.17'
Set-TextValue $ws "E25" '  -6.81%  '

# Row 26
Set-TextValue $ws "D26" '2.261'
Set-TextValue $ws "E26" '  -7.14%  '

# Row 27
Set-TextValue $ws "D27" '142.86'
Set-TextValue $ws "E27" '  -3.89%  '

# Row 28
Set-TextValue $ws "D28" '17.26'
Set-TextValue $ws "E28" '  -7.45%  '

# Row 29
Set-TextValue $ws "D29" '1.619.04'
Set-TextValue $ws "E29" '  -6.56%  '

# Row 30
Set-TextValue $ws "D30" '110.21'
Set-TextValue $ws "E30" '  -7.94%  '

# Row 31
Set-TextValue $ws "D31" '0.9245'
Set-TextValue $ws "E31" '  -6.98%  '

# Row 32
Set-TextValue $ws "B32" 'HuobiToken'
Set-TextValue $ws "C32" 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
Set-TextValue $ws "D32" '3.623'
Set-TextValue $ws "E32" '  -25.64%  '

# Row 33
Set-TextValue $ws "B33" 'Filecoin'
Set-TextValue $ws "C33" 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue $ws "D33" '5.497'
Set-TextValue $ws "E33" '  -7.01%  '

# Row 34
Set-TextValue $ws "D34" '0.07747'
Set-TextValue $ws "E34" '  -7.20%  '

# Row 35
Set-TextValue $ws "D35" '8.295'
Set-TextValue $ws "E35" '  -9.68%  '

# Row 36
Set-TextValue $ws "D36" '1.440'
Set-TextValue $ws "E36" '  -10.40%  '

# Row 37
Set-TextValue $ws "B37" 'Aptos'
Set-TextValue $ws "C37" 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextValue $ws "D37" '11.07'
Set-TextValue $ws "E37" '  +2.39%  '

# Row 38
Set-TextValue $ws "B38" 'Frax'
Set-TextValue $ws "C38" 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
Set-TextValue $ws "D38" '1.009'
Set-TextValue $ws "E38" '  +0.96%  '

# Row 39
Set-TextValue $ws "D39" '4.804'
Set-TextValue $ws "E39" '  -6.92%  '

# Row 40
Set-TextValue $ws "D40" '0.05653'
Set-TextValue $ws "E40" '  -6.35%  '

# Row 41
Set-TextValue $ws "D41" '0.02061'
Set-TextValue $ws "E41" '  -9.63%  '

# Row 42
Set-TextValue $ws "D42" '0.1924'
Set-TextValue $ws "E42" '  -6.48%  '

# Row 43
Set-TextValue $ws "D43" '1.122'
Set-TextValue $ws "E43" '  -7.35%  '

# Row 44
Set-TextValue $ws "D44" '3.594'
Set-TextValue $ws "E44" '  -4.86%  '

# Row 45
Set-TextValue $ws "D45" '0.5364'
Set-TextValue $ws "E45" '  -8.29%  '

# Row 46
Set-TextValue $ws "D46" '12.12'
Set-TextValue $ws "E46" '  -7.87%  '

# Row 47
Set-TextValue $ws "D47" '0.5186'
Set-TextValue $ws "E47" '  -7.68%  '

# Row 48
Set-TextValue $ws "D48" '1.785'
Set-TextValue $ws "E48" '  -6.41%  '

# Row 49
Set-TextValue $ws "D49" '110.80'
Set-TextValue $ws "E49" '  -6.78%  '

# Row 50
Set-TextValue $ws "D50" '1.069'
Set-TextValue $ws "E50" '  -6.48%  '

# Row 51
Set-TextValue $ws "D51" '0.06346'
Set-TextValue $ws "E51" '  -7.06%  '
